$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B/C shift right to C/D)
$ws.Columns.Item(2).Insert()

# New column B gets the same width as column A
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Header for the new StatQuery column
$ws.Range("B1").Value = "StatQuery"

# New query text for row 2, same wrap-text style as A2
$statQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Malignant neoplasm of the respiratory tract cell type specified :: Lung adenocarcinoma (metastatic)']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$ws.Range("B2").Value = $statQuery
$ws.Range("B2").WrapText = $true

# Update the active selection to B2 only
$ws.Range("B2").Select() | Out-Null
